$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix "Objetivos:" row (row 10) B/C text
$ws.Range("B10").Value = 'Informar os alunos dos produtos cerâmicos, de suas propriedades e suas aplicações'
$ws.Range("C10").Value = 'Informar os alunos dos produtos cerâmicos, de suas propriedades e suas aplicações'

# 2. Insert new row at 13 for the professor-name row (moved up from old mis-placed spot)
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
# Copy B14:C14 formats (the row that will hold the target style after shift) onto B13:C13
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B13").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C13").Value = '1922320 - Sebastiao Ribeiro'

# 3. Row 14 (was 13, "Programa resumido:") - replace B/C text
$ws.Range("B14").Value = '1.Cerâmica vermelha2.Cerâmica branca3.Refratários4.Cerâmica eletro-eletrônica5.Cerâmica avançada estrutural'
$ws.Range("C14").Value = '1.Cerâmica vermelha2.Cerâmica branca3.Refratários4.Cerâmica eletro-eletrônica5.Cerâmica avançada estrutural'

# 4. Row 16 (was 15, "Programa:") - replace B/C text
$ws.Range("B16").Value = '01 - Cerâmica vermelha02 - Cerâmica vermelha03 - Cerâmica vermelha04 - Cerâmica vermelha05 - Cerâmica branca06 - Cerâmica branca07 - Refratários08 - Refratários09 - Refratários10 - Cerâmica eletro-eletrônica11 - Cerâmica eletro-eletrônica12 - Cerâmica eletro-eletrônica13 - Cerâmica eletro-eletrônica14 - Cerâmica avançada estrutural15 - Cerâmica avançada estrutural'
$ws.Range("C16").Value = '01 - Cerâmica vermelha02 - Cerâmica vermelha03 - Cerâmica vermelha04 - Cerâmica vermelha05 - Cerâmica branca06 - Cerâmica branca07 - Refratários08 - Refratários09 - Refratários10 - Cerâmica eletro-eletrônica11 - Cerâmica eletro-eletrônica12 - Cerâmica eletro-eletrônica13 - Cerâmica eletro-eletrônica14 - Cerâmica avançada estrutural15 - Cerâmica avançada estrutural'

# 5. Row 19 (was 18, "Método:") - replace B/C text
$ws.Range("B19").Value = 'Duas provas escritas (P1 e P2), valendo de 0 (zero) a 10 (dez)'
$ws.Range("C19").Value = 'Duas provas escritas (P1 e P2), valendo de 0 (zero) a 10 (dez)'

# 6. Row 20 (was 19, "Critério:") - replace B/C text
$ws.Range("B20").Value = 'Média Parcial (MP): (P1 + P2)/2Média Parcial igual ou superior a 5: aprovação diretaMédia Parcial entre 3 e 5: recuperaçãoMédia Parcial inferior a 5: reprovação direta'
$ws.Range("C20").Value = 'Média Parcial (MP): (P1 + P2)/2Média Parcial igual ou superior a 5: aprovação diretaMédia Parcial entre 3 e 5: recuperaçãoMédia Parcial inferior a 5: reprovação direta'

# 7. Row 21 (was 20, "Norma de recuperação:") - replace B/C text
$ws.Range("B21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez)Média Final:(MP + PR)/2Média Final igual ou superior a 5 (cinco): aprovadoMédia Final inferior a 5: reprovado'
$ws.Range("C21").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez)Média Final:(MP + PR)/2Média Final igual ou superior a 5 (cinco): aprovadoMédia Final inferior a 5: reprovado'

# 8. Row 22 (was 21, "Bibliografia:") - replace B/C text
$ws.Range("B22").Value = '1.F. Singer & S. S. Singer, Cerâmica Industrial, V. 11, 19712.Salmang & Scholze, Keramik: Teil2 Keramische Werkstoffe, Springer Verlag, 19833.L. M. Levinson, Electronic Ceramics, Properties, Devices and Applications4.M. J. Hoffmann, Silicon Nitride'
$ws.Range("C22").Value = '1.F. Singer & S. S. Singer, Cerâmica Industrial, V. 11, 19712.Salmang & Scholze, Keramik: Teil2 Keramische Werkstoffe, Springer Verlag, 19833.L. M. Levinson, Electronic Ceramics, Properties, Devices and Applications4.M. J. Hoffmann, Silicon Nitride'
